# The deck's theme is switched from the green/teal "Integral" palette to the
# stock "Office Theme" palette. Drive it the way a user would in the UI:
# Design > Variants > Colors > (pick a different theme-color set), which in
# the object model is exposed per-slide as ThemeColorScheme and edits the
# twelve theme colors (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) that live in
# the slide master's theme part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Office Theme color scheme, in ThemeColorScheme.Colors() index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
